# Realestate Update resale numbers 2023-07-01 23:51
# Append a new data row (row 98) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 98

# Columns A-D hold text that looks like a date/number ("2023-07-01", "26"),
# so force a Text number format first to stop Excel from auto-converting
# them into a date serial / numeric value.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-07-01"

$ws.Cells.Item($row, 2).Value = "23:38:53"

$ws.Cells.Item($row, 3).Value = "Saturday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "26"

$ws.Cells.Item($row, 5).Value = 123584
$ws.Cells.Item($row, 6).Value = 134944
$ws.Cells.Item($row, 7).Value = 161777
$ws.Cells.Item($row, 8).Value = 131577
$ws.Cells.Item($row, 9).Value = 175696
$ws.Cells.Item($row, 10).Value = 113354
$ws.Cells.Item($row, 11).Value = 205043
$ws.Cells.Item($row, 12).Value = 222754
$ws.Cells.Item($row, 13).Value = 174590
$ws.Cells.Item($row, 14).Value = 103436
$ws.Cells.Item($row, 15).Value = 38906
$ws.Cells.Item($row, 16).Value = 32909
$ws.Cells.Item($row, 17).Value = 52001
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36279
$ws.Cells.Item($row, 20).Value = -1
